$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 53.64296733333333
$ws.Range("H2").Value = 160.928902
$ws.Range("I2").Value = 0.1405570576660657
$ws.Range("J2").Value = 0.1405570576660657
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.002279333333333333
$ws.Range("N2").Value = 0.006838
$ws.Range("Q2").Value = 0.1222702035417778
$ws.Range("R2").Value = 1.100431831876
$ws.Range("S2").Value = 0.1405570576660657
$ws.Range("T2").Value = 0.1405570576660657

# Row 3
$ws.Range("I3").Value = 0.83973167405618
$ws.Range("J3").Value = 0.8397316740561799
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.002279333333333333
$ws.Range("N3").Value = 0.006838
$ws.Range("Q3").Value = 0.7304803075151111
$ws.Range("R3").Value = 6.574322767636
$ws.Range("S3").Value = 0.83973167405618
$ws.Range("T3").Value = 0.8397316740561799

# Row 4
$ws.Range("G4").Value = 7.522716666666668
$ws.Range("I4").Value = 0.01971126827775425
$ws.Range("J4").Value = 0.01971126827775425
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.002279333333333333
$ws.Range("N4").Value = 0.006838
$ws.Range("Q4").Value = 0.01714677885555556
$ws.Range("R4").Value = 0.1543210097
$ws.Range("S4").Value = 0.01971126827775425
$ws.Range("T4").Value = 0.01971126827775425
